# Add a new hike entry "Mather Memorial Parkway" to the "Hike Difficulties"
# table, inserted in its sorted position as row 42 (pushing the existing
# rows 42-78 down to 43-79).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 42, shifting existing data (and the rows below)
# down by one.
$ws.Rows(42).Insert()

# Fill in the new row's data.
$ws.Range("A42").Value = "Mather Memorial Parkway"
$ws.Range("B42").Value = 10.8
$ws.Range("C42").Value = 1050
$ws.Range("D42").Value = "moderate"

# Grow Table1 so the new row (and the row that used to be last) are both
# included in the table range / autofilter / sort state.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D79"))

# Match the saved selection state from the authored edit.
$ws.Range("D42").Select()
